$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column A values for rows 2-6 first (existing rows get new text, new rows added)
$ws.Range("A2").Value = "dolore quis consectetur"
$ws.Range("A3").Value = "magna in voluptate quis"
$ws.Range("A4").Value = "adipisicing eiusmod magna in Excepteur"
$ws.Range("A5").Value = "fugiat"
$ws.Range("A6").Value = "in velit nostrud Excepteur"

# Then set column B values for rows 2-6
$ws.Range("B2").Value = "voluptate dolor tempor"
$ws.Range("B3").Value = "quis do ad velit aute"
$ws.Range("B4").Value = "culpa minim"
$ws.Range("B5").Value = "eiusmod"
$ws.Range("B6").Value = "anim Ut do"

# Then set column C values
$ws.Range("C2").Value = -1650109834
$ws.Range("C3").Value = -504719199
$ws.Range("C4").Value = -572649954
$ws.Range("C5").Value = -921315400
$ws.Range("C6").Value = -1400435088

# Then set column D values
$ws.Range("D2").Value = -1325551028
$ws.Range("D3").Value = 1835153350
$ws.Range("D4").Value = -478153369
$ws.Range("D5").Value = 1150486065
$ws.Range("D6").Value = -261493675
